$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade rows to append (row 10, 11, 12)
$rows = @(
    @{ A = 9889.33;            B = 9895.27; C = 283.47000000000003; D = 283.29000000000002; E = $false; F = -0.06;                  G = 42613.765601851854; H = $false },
    @{ A = 9836.92;            B = 9889.33; C = 282.39;              D = 280.89;              E = $false; F = -0.53;                  G = 42614.673020833332; H = $false },
    @{ A = 9865.4500000000007; B = 9836.92; C = 280.62;              D = 281.44;              E = $false; F = 0.28999999999999998;    G = 42615.750208333331; H = $true  }
)

$startRow = 10
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F

    $gCell = $ws.Cells.Item($r, 7)
    $ws.Cells.Item($r - 1, 7).Copy()
    $gCell.PasteSpecial(-4122)
    $gCell.Value = $data.G

    $ws.Cells.Item($r, 8).Value = $data.H
}
